# Auto-generated PowerShell Excel COM-interop script
# Applies updated market-price values across the Kraken_Profits sheets
# per the scheduled runner's data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 54
$ws.Range("H54").Value = '6416.3335'
$ws.Range("I54").Value = '5749.5'
$ws.Range("K54").Value = '5749.5'
$ws.Range("M54").Value = '-5263.5'
# Row 75
$ws.Range("H75").Value = '0'
$ws.Range("J75").Value = '0'
$ws.Range("L75").Value = '0'
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = '0'
$ws.Range("J78").Value = '0'
$ws.Range("L78").Value = '0'
$ws.Range("N78").ClearContents()
# Row 80
$ws.Range("H80").Value = '630'
$ws.Range("J80").Value = '630'
$ws.Range("L80").Value = '1890'
$ws.Range("N80").Value = '-3886'
# Row 83
$ws.Range("H83").Value = '630'
$ws.Range("J83").Value = '630'
$ws.Range("L83").Value = '5670'
$ws.Range("N83").Value = '-15654'
# Row 93
$ws.Range("H93").Value = '20875'
$ws.Range("J93").Value = '20875'
$ws.Range("L93").Value = '20875'
$ws.Range("N93").Value = '-25867'
# Row 99
$ws.Range("H99").Value = '292'
$ws.Range("I99").Value = '315'
$ws.Range("K99").Value = '945'
$ws.Range("M99").Value = '553'
# Row 103
$ws.Range("H103").Value = '1494'
$ws.Range("I103").Value = '0'
$ws.Range("J103").Value = '1494'
$ws.Range("K103").Value = '0'
$ws.Range("L103").Value = '4482'
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = '-5654'
# Row 105
$ws.Range("H105").Value = '30250'
$ws.Range("J105").Value = '30250'
$ws.Range("L105").Value = '30250'
$ws.Range("N105").Value = '-37238'
# Row 110
$ws.Range("H110").Value = '10000'
$ws.Range("J110").Value = '10000'
$ws.Range("L110").Value = '10000'
$ws.Range("N110").Value = '-18180'
# Row 111
$ws.Range("H111").Value = '331.66666'
$ws.Range("I111").Value = '335'
$ws.Range("K111").Value = '1005'
$ws.Range("M111").Value = '2062'
# Row 113
$ws.Range("H113").Value = '3249.75'
$ws.Range("I113").Value = '3499.5'
$ws.Range("K113").Value = '3499.5'
$ws.Range("M113").Value = '-245.5'
# Row 117
$ws.Range("H117").Value = '50000'
$ws.Range("J117").Value = '50000'
$ws.Range("L117").Value = '50000'
$ws.Range("N117").Value = '-59178'
# Row 118
$ws.Range("H118").Value = '1200'
$ws.Range("I118").Value = '1200'
$ws.Range("K118").Value = '3600'
$ws.Range("M118").Value = '-1943'
# Row 132
$ws.Range("H132").Value = '7387.7'
$ws.Range("I132").Value = '7315.1665'
$ws.Range("K132").Value = '21945.4995'
$ws.Range("M132").Value = '-19415.4995'

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = '3220.2222'
$ws.Range("I97").Value = '2580.3333'
$ws.Range("J97").Value = '4500'
$ws.Range("K97").Value = '2580.3333'
$ws.Range("L97").Value = '4500'
$ws.Range("M97").Value = '-2084.3333'
$ws.Range("N97").Value = '-5492'
# Row 125
$ws.Range("H125").Value = '12756429'
$ws.Range("J125").Value = '12756429'
$ws.Range("L125").Value = '12756429'
$ws.Range("N125").Value = '-12766269'

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = '5099.8'
$ws.Range("J105").Value = '6500'
$ws.Range("L105").Value = '6500'
$ws.Range("N105").Value = '-9994'

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = '1233.5'
$ws.Range("I62").Value = '1234'
$ws.Range("J62").Value = '1233'
$ws.Range("K62").Value = '1234'
$ws.Range("L62").Value = '1233'
$ws.Range("M62").Value = '-610'
$ws.Range("N62").Value = '-2481'
# Row 65
$ws.Range("H65").Value = '1233.5'
$ws.Range("I65").Value = '1234'
$ws.Range("J65").Value = '1233'
$ws.Range("K65").Value = '6170'
$ws.Range("L65").Value = '6165'
$ws.Range("M65").Value = '-3050'
$ws.Range("N65").Value = '-12405'
# Row 99
$ws.Range("H99").Value = '1800'
$ws.Range("I99").Value = '0'
$ws.Range("J99").Value = '1800'
$ws.Range("K99").Value = '0'
$ws.Range("L99").Value = '1800'
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = '-4796'
# Row 107
$ws.Range("H107").Value = '773.75'
$ws.Range("I107").Value = '773.75'
$ws.Range("K107").Value = '773.75'
$ws.Range("M107").Value = '1146.25'
# Row 122
$ws.Range("H122").Value = '1779.8'
$ws.Range("I122").Value = '1724.75'
$ws.Range("K122").Value = '5174.25'
$ws.Range("M122").Value = '-2724.25'
# Row 124
$ws.Range("H124").Value = '50326'
$ws.Range("J124").Value = '50326'
$ws.Range("L124").Value = '50326'
$ws.Range("N124").Value = '-55236'
# Row 126
$ws.Range("H126").Value = '1800'
$ws.Range("I126").Value = '0'
$ws.Range("J126").Value = '1800'
$ws.Range("K126").Value = '0'
$ws.Range("L126").Value = '5400'
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = '-10340'

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = '1595'
$ws.Range("I75").Value = '0'
$ws.Range("J75").Value = '1595'
$ws.Range("K75").Value = '0'
$ws.Range("L75").Value = '4785'
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = '-6781'
# Row 78
$ws.Range("H78").Value = '1595'
$ws.Range("I78").Value = '0'
$ws.Range("J78").Value = '1595'
$ws.Range("K78").Value = '0'
$ws.Range("L78").Value = '14355'
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = '-24339'

$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = '102.5'
$ws.Range("I19").Value = '102.5'
$ws.Range("K19").Value = '102.5'
$ws.Range("M19").Value = '185.5'
# Row 80
$ws.Range("H80").Value = '3894.1667'
$ws.Range("I80").Value = '3273'
$ws.Range("J80").Value = '7000'
$ws.Range("K80").Value = '3273'
$ws.Range("L80").Value = '7000'
$ws.Range("M80").Value = '-2275'
$ws.Range("N80").Value = '-8996'
# Row 83
$ws.Range("H83").Value = '3894.1667'
$ws.Range("I83").Value = '3273'
$ws.Range("J83").Value = '7000'
$ws.Range("K83").Value = '16365'
$ws.Range("L83").Value = '35000'
$ws.Range("M83").Value = '-11373'
$ws.Range("N83").Value = '-44984'

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = '2018.8572'
$ws.Range("J16").Value = '1230'
$ws.Range("L16").Value = '1230'
$ws.Range("N16").Value = '-1570'
# Row 55
$ws.Range("H55").Value = '2075.5'
$ws.Range("I55").Value = '301'
$ws.Range("J55").Value = '2667'
$ws.Range("K55").Value = '301'
$ws.Range("L55").Value = '2667'
$ws.Range("M55").Value = '-128'
$ws.Range("N55").Value = '-3013'
# Row 93
$ws.Range("H93").Value = '8812.143'
$ws.Range("I93").Value = '8812.143'
$ws.Range("K93").Value = '8812.143'
$ws.Range("M93").Value = '-7564.143'
# Row 122
$ws.Range("H122").Value = '4324.875'
$ws.Range("I122").Value = '4371.2856'
$ws.Range("K122").Value = '13113.8568'
$ws.Range("M122").Value = '-10663.8568'
# Row 132
$ws.Range("H132").Value = '4020'
$ws.Range("I132").Value = '4020'
$ws.Range("K132").Value = '12060'
$ws.Range("M132").Value = '-9530'

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = '1736.5'
$ws.Range("I126").Value = '1736.5'
$ws.Range("K126").Value = '5209.5'
$ws.Range("M126").Value = '-2739.5'

